# Adding CMAD presentation -ver2
# Appends two new "Title and Content" slides to the end of the deck:
#   12: React JS Components
#   13: React Component Classes
#
# Slide IDs are assigned by PowerPoint in creation order, so to reproduce the
# authored id ordering (268 appears before 267 in the sldIdLst) the slide
# that is inserted SECOND (and therefore receives the higher internal id)
# must be placed BEFORE the first one in the final slide order. We do this
# by first creating the "React Component Classes" slide at the insertion
# point (temporarily at index 12) and then inserting the "React JS
# Components" slide in front of it at the same index 12, which pushes the
# first slide down to index 13.

$p = $ppt.ActivePresentation
$insertAt = $p.Slides.Count + 1

# Step 1: create the slide that will end up SECOND ("React Component Classes").
$sClasses = $p.Slides.Add($insertAt, 2)
$sClasses.Shapes.Item(1).TextFrame.TextRange.Text = "React Component Classes"

$tfClasses = $sClasses.Shapes.Item(2).TextFrame
$tfClasses.TextRange.Text = "Class CMADApp extends React.Component{`r      constructor(props){`r        super(props)`r        this.messages = <List of messages>`r      }`r      Class Header extends React.Component{  #“static header”}`r      Class MessageStats extends React.Component{#“Display based on Stats”}`r      Class MessageFilter extends React.Component{ #“Display based on filters”}`r      Class MessageScroll extends React.Component { #“Scrolling live display”}`r}`rThe messages will be state which will be passed as props into the MessageStats and MessageScroll. We need to periodically poll the server for message updates."
$tfClasses.TextRange.ParagraphFormat.Bullet.Visible = 0
$tfClasses.AutoSize = 2

# Step 2: create the slide that will end up FIRST ("React JS Components"),
# inserting it at the same index so it lands right before the slide above.
$sComponents = $p.Slides.Add($insertAt, 2)
$sComponents.Shapes.Item(1).TextFrame.TextRange.Text = "React JS Components"

$tfComponents = $sComponents.Shapes.Item(2).TextFrame
$tfComponents.TextRange.Text = "Class App extends React.Component {`rRender()`r   return`r      <div> `r           <Header />`r`t< MessageStats />`r    `t< MessageFilter />`r`t<MessageScroll />`r     </div>`r}`r         "
$tfComponents.TextRange.ParagraphFormat.Bullet.Visible = 0
$tfComponents.AutoSize = 2
